$d = $word.ActiveDocument
$nbsp = [string]([char]0x00A0)

# Locate the first "wordpress" occurrence together with the non-breaking
# space that follows it (so the run holding "wordpress" - and the
# w:proofErr spell-check markers wrapping it - gets fully replaced instead
# of leaving stray markers behind), then retype it as "blogger" followed
# by that same non-breaking space.
$rWord = $d.Content
$rWord.Find.Execute("wordpress" + $nbsp, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$wpStart = $rWord.Start
$rWord.Text = "blogger" + $nbsp

# Re-derive the text once more, now that "wordpress" became "blogger", to
# recompute the run split points we still need.
$full = $d.Content.Text
$bloggerEnd = $wpStart + 7          # end of the word "blogger"
$travaillerStart = $full.IndexOf("Travailler sur toutes")

# Re-create the original run boundary right after "blogger" (before the
# " ? " that follows) and right before "Travailler sur toutes les données"
# by dropping a temporary bookmark at each spot (which forces a run split)
# and immediately deleting it again (the split itself persists).
foreach ($pos in @($bloggerEnd, $travaillerStart)) {
    $rs = $d.Range($pos, $pos)
    $d.Bookmarks.Add("TmpSplit", $rs)
    $d.Bookmarks("TmpSplit").Delete()
}

# Finally, stamp the _GoBack bookmark at the spot of the last real edit:
# between "uniquem" and "ent" of "uniquement" (Word tracks the caret
# position of the latest change there). Adding a bookmark under a name
# that already exists relocates it, so the stale one that used to sit
# alone in the trailing empty paragraph disappears automatically.
$rSplit = $d.Content
$rSplit.Find.Execute("uniquem", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rSplit.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rSplit)
